$d = $word.ActiveDocument

# Locate the run containing "GNS Healthcare, Cambridge MA," and split it into
# three runs: "GNS Healthcare" / " (Rebranded Aitia Bio Jan 2023)" / ", Cambridge MA,"
# all sharing identical formatting (Tahoma, smallCaps, sz 20 == 10pt).

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("GNS Healthcare", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$ins = $find.Parent.Duplicate
$ins.Collapse(0)
$ins.InsertAfter(" (Rebranded Aitia Bio Jan 2023)")

# The newly inserted text currently shares formatting with its neighbours, so the
# engine coalesces it back into one run. Briefly nudge its size to force the split
# into a separate run, then restore the original size (10pt / w:sz=20) via a fresh
# Find so the final formatting is byte-identical to its neighbours.
$findNew = $d.Content.Find
$findNew.ClearFormatting()
$foundNew = $findNew.Execute(" (Rebranded Aitia Bio Jan 2023)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findNew.Parent.Font.Size = 11

$findNew2 = $d.Content.Find
$findNew2.ClearFormatting()
$foundNew2 = $findNew2.Execute(" (Rebranded Aitia Bio Jan 2023)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findNew2.Parent.Font.Size = 10
